# Generate Report for Handoff
# Updates the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps for the
# file 0f8c911a-4486-4231-803f-7dce16c12ccb.md (row 5 on every sheet) after a new
# handoff xliff was generated for it.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G5").Value = "2016-08-17 22:40:56"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn.Range("H5").Value = "2016-08-17 22:40:51"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe.Range("H5").Value = "2016-08-17 22:40:56"
